# Daily attendance processing - 2026-01-10 04:25:08
# Normalize the "Recorded By" column (G) so that "System" is always listed
# first among the recorder names, preserving the relative order of the
# remaining entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $raw = $cell.Value()

    if ($raw -eq $null) {
        continue
    }

    $text = [string]$raw
    if ($text -eq "") {
        continue
    }

    $parts = $text -split ","
    for ($i = 0; $i -lt $parts.Length; $i++) {
        $parts[$i] = $parts[$i].Trim()
    }

    $systemIndex = -1
    for ($i = 0; $i -lt $parts.Length; $i++) {
        if ($parts[$i] -eq "System") {
            $systemIndex = $i
            break
        }
    }

    if ($systemIndex -gt 0) {
        $newParts = @("System")
        for ($i = 0; $i -lt $parts.Length; $i++) {
            if ($i -ne $systemIndex) {
                $newParts += $parts[$i]
            }
        }
        $newText = [string]::Join(", ", $newParts)
        $cell.Value = $newText
    }
}
